$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Patient Selections")
$ws2 = $wb.Worksheets.Item("Health-ITUES")

# Add question header labels Q1..Q17 to Health-ITUES row 1, columns C..S
$headers = @("Q1","Q2","Q3","Q4","Q5","Q6","Q7","Q8","Q9","Q10","Q11","Q12","Q13","Q14","Q15","Q16","Q17")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 3 + $i
    $ws2.Cells.Item(1, $col).Value = $headers[$i]
}

# Update selections / active sheet state
$ws2.Range("F22").Select()
$ws1.Activate()
$ws1.Range("D6").Select()
